$d = $word.ActiveDocument

# 1. Remove the "Meta description" paragraph that follows the H1 heading
#    ("Play Geisha for Free - A Japanese-Themed Slot Game").
$d.Paragraphs.Item(2).Range.Delete()

# 2. Insert a new bold "Play Geisha for Free - A Japanese-Themed Slot Game"
#    paragraph right before the final "Prompt: ..." paragraph.
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$lastPara.Range.InsertParagraphBefore()
$newPara = $d.Paragraphs.Item($count)
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Geisha for Free - A Japanese-Themed Slot Game</w:t></w:r></w:p>'
$newPara.Range.InsertXML($xml)

# 3. Replace the text of the (now last) "Prompt: ..." paragraph with the new
#    meta-description copy, keeping its italic formatting intact.
$finalPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$finalPara.Range.Find.Execute("Prompt: Create a feature image for Aristocrat's Geisha slot game in cartoon style depicting a happy Maya warrior with glasses. For this feature image, the artist could draw a Maya warrior in a cartoon style with a smile on their face and glasses on their eyes. The warrior could be holding a Geisha girl's hand, and they could both be standing in front of a background with a Mount Fuji and a dragon. The artist could also include symbols from the game, such as the lotus flower and playing cards, to help tie the image back to the slot game. The overall vibe of the image should be bright, cheerful, and playful to capture the fun and entertaining nature of the game.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Discover Geisha, a Japanese-themed slot game from Aristocrat. Play for free and explore captivating graphics, engaging gameplay, and wild and scatter symbols for big wins.", 2)
